# Round the computed ConvexHullArea figures (column D) to whole numbers,
# matching the "Add separability results in csv" cleanup.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 202

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -ne $null) {
        $val = [double]$val
        if ($val -ge 0) {
            $rounded = [Math]::Floor($val + 0.5)
        } else {
            $rounded = -([Math]::Floor(-$val + 0.5))
        }
        $cell.Value = $rounded
    }
}
